$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "Hall Short Name" ---------------------------------
# Header cell E4 mirrors the formatting of the existing header row
# (A4:C4 use style index 2), so copy format from A4 then overwrite value.
$ws.Range("A4").Copy()
$ws.Range("E4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E4").Value = "Hall Short Name"

# Existing students (rows 5-8) all belong to the "RTH" hall.
$ws.Range("E5").Value = "RTH"
$ws.Range("E6").Value = "RTH"
$ws.Range("E7").Value = "RTH"
$ws.Range("E8").Value = "RTH"

# --- New student record (row 9) ---------------------------------------
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Afrin Ahmed Eva"
$ws.Range("C9").Value = 41140
$ws.Range("D9").Value = "F"
$ws.Range("E9").Value = "SHH"

# Column E width to comfortably fit the new header text.
$ws.Columns("E").ColumnWidth = 20.6

# Match the workbook's last active selection after the edit.
$ws.Range("E8").Select()
